$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 new rows starting at row 144, pushing existing rows 144-170 down to 157-183
$ws.Range("A144:A156").EntireRow.Insert()

# Ensure column A keeps the bold/bordered/centered style used throughout the data rows
$colA = $ws.Range("A144:A183")
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = 1

# Row 144
$ws.Cells.Item(144, 1).Value = 142
$ws.Cells.Item(144, 2).Value = 3
$ws.Cells.Item(144, 3).Value = "What are the areas within a four-minute drive of each fire station at 2 a.m. on Tuesday in Utrecht"
$ws.Cells.Item(144, 4).Value = "Utrecht"
$ws.Cells.Item(144, 5).Value = ""
$ws.Cells.Item(144, 6).Value = ""
$ws.Cells.Item(144, 7).Value = ""
$ws.Cells.Item(144, 8).Value = "done"
$ws.Cells.Item(144, 9).Value = "Data queries"
$ws.Cells.Item(144, 10).Value = "network analysis"
$ws.Cells.Item(144, 11).Value = "classification"
$ws.Cells.Item(144, 12).Value = "Data queries"
$ws.Cells.Item(144, 13).Value = "Overlay analysis"
$ws.Cells.Item(144, 14).Value = ""
$ws.Cells.Item(144, 15).Value = ""
$ws.Cells.Item(144, 16).Value = ""
$ws.Cells.Item(144, 17).Value = ""
$ws.Cells.Item(144, 18).Value = ""
$ws.Cells.Item(144, 19).Value = ""
$ws.Cells.Item(144, 20).Value = ""
$ws.Cells.Item(144, 21).Value = ""
$ws.Cells.Item(144, 22).Value = ""
$ws.Cells.Item(144, 23).Value = ""
$ws.Cells.Item(144, 24).Value = ""
$ws.Cells.Item(144, 25).Value = "data queries,network analysis,classification,data queries,overlay analysis"
$ws.Cells.Item(144, 26).Value = 0
$ws.Cells.Item(144, 27).Value = $false

# Row 145
$ws.Cells.Item(145, 1).Value = 143
$ws.Cells.Item(145, 2).Value = 13
$ws.Cells.Item(145, 3).Value = "What area are within 50 km from family physician services in Saskatchewan in Canada"
$ws.Cells.Item(145, 4).Value = "Saskatchewan"
$ws.Cells.Item(145, 5).Value = "Canada"
$ws.Cells.Item(145, 6).Value = ""
$ws.Cells.Item(145, 7).Value = "amenity=doctor"
$ws.Cells.Item(145, 8).Value = "done"
$ws.Cells.Item(145, 9).Value = "Data queries"
$ws.Cells.Item(145, 10).Value = "Buffer"
$ws.Cells.Item(145, 11).Value = "Overlay analysis"
$ws.Cells.Item(145, 12).Value = ""
$ws.Cells.Item(145, 13).Value = ""
$ws.Cells.Item(145, 14).Value = ""
$ws.Cells.Item(145, 15).Value = ""
$ws.Cells.Item(145, 16).Value = ""
$ws.Cells.Item(145, 17).Value = ""
$ws.Cells.Item(145, 18).Value = ""
$ws.Cells.Item(145, 19).Value = ""
$ws.Cells.Item(145, 20).Value = ""
$ws.Cells.Item(145, 21).Value = ""
$ws.Cells.Item(145, 22).Value = ""
$ws.Cells.Item(145, 23).Value = ""
$ws.Cells.Item(145, 24).Value = ""
$ws.Cells.Item(145, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(145, 26).Value = 10
$ws.Cells.Item(145, 27).Value = $false

# Row 146
$ws.Cells.Item(146, 1).Value = 144
$ws.Cells.Item(146, 2).Value = 16
$ws.Cells.Item(146, 3).Value = "What areas are inside 1000 foot of schools in El Cajon"
$ws.Cells.Item(146, 4).Value = "El Cajon"
$ws.Cells.Item(146, 5).Value = ""
$ws.Cells.Item(146, 6).Value = ""
$ws.Cells.Item(146, 7).Value = "amenity=school, amenity=kindergarten"
$ws.Cells.Item(146, 8).Value = "done"
$ws.Cells.Item(146, 9).Value = "Data queries"
$ws.Cells.Item(146, 10).Value = "Buffer"
$ws.Cells.Item(146, 11).Value = "Overlay analysis"
$ws.Cells.Item(146, 12).Value = ""
$ws.Cells.Item(146, 13).Value = ""
$ws.Cells.Item(146, 14).Value = ""
$ws.Cells.Item(146, 15).Value = ""
$ws.Cells.Item(146, 16).Value = ""
$ws.Cells.Item(146, 17).Value = ""
$ws.Cells.Item(146, 18).Value = ""
$ws.Cells.Item(146, 19).Value = ""
$ws.Cells.Item(146, 20).Value = ""
$ws.Cells.Item(146, 21).Value = ""
$ws.Cells.Item(146, 22).Value = ""
$ws.Cells.Item(146, 23).Value = ""
$ws.Cells.Item(146, 24).Value = ""
$ws.Cells.Item(146, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(146, 26).Value = 10
$ws.Cells.Item(146, 27).Value = $false

# Row 147
$ws.Cells.Item(147, 1).Value = 145
$ws.Cells.Item(147, 2).Value = 20
$ws.Cells.Item(147, 3).Value = "What areas are green belt areas in Houston"
$ws.Cells.Item(147, 4).Value = "Houston"
$ws.Cells.Item(147, 5).Value = ""
$ws.Cells.Item(147, 6).Value = ""
$ws.Cells.Item(147, 7).Value = ""
$ws.Cells.Item(147, 8).Value = "done"
$ws.Cells.Item(147, 9).Value = "Data queries"
$ws.Cells.Item(147, 10).Value = "Overlay analysis"
$ws.Cells.Item(147, 11).Value = ""
$ws.Cells.Item(147, 12).Value = ""
$ws.Cells.Item(147, 13).Value = ""
$ws.Cells.Item(147, 14).Value = ""
$ws.Cells.Item(147, 15).Value = ""
$ws.Cells.Item(147, 16).Value = ""
$ws.Cells.Item(147, 17).Value = ""
$ws.Cells.Item(147, 18).Value = ""
$ws.Cells.Item(147, 19).Value = ""
$ws.Cells.Item(147, 20).Value = ""
$ws.Cells.Item(147, 21).Value = ""
$ws.Cells.Item(147, 22).Value = ""
$ws.Cells.Item(147, 23).Value = ""
$ws.Cells.Item(147, 24).Value = ""
$ws.Cells.Item(147, 25).Value = "data queries,overlay analysis"
$ws.Cells.Item(147, 26).Value = 13
$ws.Cells.Item(147, 27).Value = $false

# Row 148
$ws.Cells.Item(148, 1).Value = 146
$ws.Cells.Item(148, 2).Value = 23
$ws.Cells.Item(148, 3).Value = "What areas are not conatined as green belt areas in Houston"
$ws.Cells.Item(148, 4).Value = "Houston"
$ws.Cells.Item(148, 5).Value = ""
$ws.Cells.Item(148, 6).Value = ""
$ws.Cells.Item(148, 7).Value = ""
$ws.Cells.Item(148, 8).Value = "done"
$ws.Cells.Item(148, 9).Value = "Data queries"
$ws.Cells.Item(148, 10).Value = "Overlay analysis"
$ws.Cells.Item(148, 11).Value = ""
$ws.Cells.Item(148, 12).Value = ""
$ws.Cells.Item(148, 13).Value = ""
$ws.Cells.Item(148, 14).Value = ""
$ws.Cells.Item(148, 15).Value = ""
$ws.Cells.Item(148, 16).Value = ""
$ws.Cells.Item(148, 17).Value = ""
$ws.Cells.Item(148, 18).Value = ""
$ws.Cells.Item(148, 19).Value = ""
$ws.Cells.Item(148, 20).Value = ""
$ws.Cells.Item(148, 21).Value = ""
$ws.Cells.Item(148, 22).Value = ""
$ws.Cells.Item(148, 23).Value = ""
$ws.Cells.Item(148, 24).Value = ""
$ws.Cells.Item(148, 25).Value = "data queries,overlay analysis"
$ws.Cells.Item(148, 26).Value = 13
$ws.Cells.Item(148, 27).Value = $false

# Row 149
$ws.Cells.Item(149, 1).Value = 147
$ws.Cells.Item(149, 2).Value = 25
$ws.Cells.Item(149, 3).Value = "What areas are not park in Houston"
$ws.Cells.Item(149, 4).Value = "Houston"
$ws.Cells.Item(149, 5).Value = ""
$ws.Cells.Item(149, 6).Value = ""
$ws.Cells.Item(149, 7).Value = "leisure=park"
$ws.Cells.Item(149, 8).Value = "done"
$ws.Cells.Item(149, 9).Value = "Data queries"
$ws.Cells.Item(149, 10).Value = ""
$ws.Cells.Item(149, 11).Value = ""
$ws.Cells.Item(149, 12).Value = ""
$ws.Cells.Item(149, 13).Value = ""
$ws.Cells.Item(149, 14).Value = ""
$ws.Cells.Item(149, 15).Value = ""
$ws.Cells.Item(149, 16).Value = ""
$ws.Cells.Item(149, 17).Value = ""
$ws.Cells.Item(149, 18).Value = ""
$ws.Cells.Item(149, 19).Value = ""
$ws.Cells.Item(149, 20).Value = ""
$ws.Cells.Item(149, 21).Value = ""
$ws.Cells.Item(149, 22).Value = ""
$ws.Cells.Item(149, 23).Value = ""
$ws.Cells.Item(149, 24).Value = ""
$ws.Cells.Item(149, 25).Value = "data queries"
$ws.Cells.Item(149, 26).Value = 1
$ws.Cells.Item(149, 27).Value = $false

# Row 150
$ws.Cells.Item(150, 1).Value = 148
$ws.Cells.Item(150, 2).Value = 26
$ws.Cells.Item(150, 3).Value = "What areas are not wetlands in Houston"
$ws.Cells.Item(150, 4).Value = "Houston"
$ws.Cells.Item(150, 5).Value = ""
$ws.Cells.Item(150, 6).Value = ""
$ws.Cells.Item(150, 7).Value = "natural=wetland"
$ws.Cells.Item(150, 8).Value = "done"
$ws.Cells.Item(150, 9).Value = "Data queries"
$ws.Cells.Item(150, 10).Value = ""
$ws.Cells.Item(150, 11).Value = ""
$ws.Cells.Item(150, 12).Value = ""
$ws.Cells.Item(150, 13).Value = ""
$ws.Cells.Item(150, 14).Value = ""
$ws.Cells.Item(150, 15).Value = ""
$ws.Cells.Item(150, 16).Value = ""
$ws.Cells.Item(150, 17).Value = ""
$ws.Cells.Item(150, 18).Value = ""
$ws.Cells.Item(150, 19).Value = ""
$ws.Cells.Item(150, 20).Value = ""
$ws.Cells.Item(150, 21).Value = ""
$ws.Cells.Item(150, 22).Value = ""
$ws.Cells.Item(150, 23).Value = ""
$ws.Cells.Item(150, 24).Value = ""
$ws.Cells.Item(150, 25).Value = "data queries"
$ws.Cells.Item(150, 26).Value = 1
$ws.Cells.Item(150, 27).Value = $false

# Row 151
$ws.Cells.Item(151, 1).Value = 149
$ws.Cells.Item(151, 2).Value = 30
$ws.Cells.Item(151, 3).Value = "What areas are outside 250 meters of human settlement in the Cape Peninsula"
$ws.Cells.Item(151, 4).Value = "the Cape Peninsula"
$ws.Cells.Item(151, 5).Value = ""
$ws.Cells.Item(151, 6).Value = ""
$ws.Cells.Item(151, 7).Value = "residential=*"
$ws.Cells.Item(151, 8).Value = "done"
$ws.Cells.Item(151, 9).Value = "Data queries"
$ws.Cells.Item(151, 10).Value = "Buffer"
$ws.Cells.Item(151, 11).Value = "Overlay analysis"
$ws.Cells.Item(151, 12).Value = ""
$ws.Cells.Item(151, 13).Value = ""
$ws.Cells.Item(151, 14).Value = ""
$ws.Cells.Item(151, 15).Value = ""
$ws.Cells.Item(151, 16).Value = ""
$ws.Cells.Item(151, 17).Value = ""
$ws.Cells.Item(151, 18).Value = ""
$ws.Cells.Item(151, 19).Value = ""
$ws.Cells.Item(151, 20).Value = ""
$ws.Cells.Item(151, 21).Value = ""
$ws.Cells.Item(151, 22).Value = ""
$ws.Cells.Item(151, 23).Value = ""
$ws.Cells.Item(151, 24).Value = ""
$ws.Cells.Item(151, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(151, 26).Value = 10
$ws.Cells.Item(151, 27).Value = $false

# Row 152
$ws.Cells.Item(152, 1).Value = 150
$ws.Cells.Item(152, 2).Value = 32
$ws.Cells.Item(152, 3).Value = "What areas are outside 3000 meters of the rivers in Spain"
$ws.Cells.Item(152, 4).Value = "Spain"
$ws.Cells.Item(152, 5).Value = ""
$ws.Cells.Item(152, 6).Value = ""
$ws.Cells.Item(152, 7).Value = "waterway=river"
$ws.Cells.Item(152, 8).Value = "done"
$ws.Cells.Item(152, 9).Value = "Data queries"
$ws.Cells.Item(152, 10).Value = "Buffer"
$ws.Cells.Item(152, 11).Value = "Overlay analysis"
$ws.Cells.Item(152, 12).Value = ""
$ws.Cells.Item(152, 13).Value = ""
$ws.Cells.Item(152, 14).Value = ""
$ws.Cells.Item(152, 15).Value = ""
$ws.Cells.Item(152, 16).Value = ""
$ws.Cells.Item(152, 17).Value = ""
$ws.Cells.Item(152, 18).Value = ""
$ws.Cells.Item(152, 19).Value = ""
$ws.Cells.Item(152, 20).Value = ""
$ws.Cells.Item(152, 21).Value = ""
$ws.Cells.Item(152, 22).Value = ""
$ws.Cells.Item(152, 23).Value = ""
$ws.Cells.Item(152, 24).Value = ""
$ws.Cells.Item(152, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(152, 26).Value = 10
$ws.Cells.Item(152, 27).Value = $false

# Row 153
$ws.Cells.Item(153, 1).Value = 151
$ws.Cells.Item(153, 2).Value = 37
$ws.Cells.Item(153, 3).Value = "What areas are within 10 miles of current transmission lines with a voltage greater than 400 in Colorado"
$ws.Cells.Item(153, 4).Value = "Colorado"
$ws.Cells.Item(153, 5).Value = ""
$ws.Cells.Item(153, 6).Value = ""
$ws.Cells.Item(153, 7).Value = "power=line"
$ws.Cells.Item(153, 8).Value = "done"
$ws.Cells.Item(153, 9).Value = "Data queries"
$ws.Cells.Item(153, 10).Value = "Buffer"
$ws.Cells.Item(153, 11).Value = "Overlay analysis"
$ws.Cells.Item(153, 12).Value = ""
$ws.Cells.Item(153, 13).Value = ""
$ws.Cells.Item(153, 14).Value = ""
$ws.Cells.Item(153, 15).Value = ""
$ws.Cells.Item(153, 16).Value = ""
$ws.Cells.Item(153, 17).Value = ""
$ws.Cells.Item(153, 18).Value = ""
$ws.Cells.Item(153, 19).Value = ""
$ws.Cells.Item(153, 20).Value = ""
$ws.Cells.Item(153, 21).Value = ""
$ws.Cells.Item(153, 22).Value = ""
$ws.Cells.Item(153, 23).Value = ""
$ws.Cells.Item(153, 24).Value = ""
$ws.Cells.Item(153, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(153, 26).Value = 10
$ws.Cells.Item(153, 27).Value = $false

# Row 154
$ws.Cells.Item(154, 1).Value = 152
$ws.Cells.Item(154, 2).Value = 47
$ws.Cells.Item(154, 3).Value = "What areas are within 2000 meters of the playgrounds in Oleander"
$ws.Cells.Item(154, 4).Value = "Oleander"
$ws.Cells.Item(154, 5).Value = ""
$ws.Cells.Item(154, 6).Value = ""
$ws.Cells.Item(154, 7).Value = "leisure=playground"
$ws.Cells.Item(154, 8).Value = "done"
$ws.Cells.Item(154, 9).Value = "Data queries"
$ws.Cells.Item(154, 10).Value = "Buffer"
$ws.Cells.Item(154, 11).Value = "Overlay analysis"
$ws.Cells.Item(154, 12).Value = ""
$ws.Cells.Item(154, 13).Value = ""
$ws.Cells.Item(154, 14).Value = ""
$ws.Cells.Item(154, 15).Value = ""
$ws.Cells.Item(154, 16).Value = ""
$ws.Cells.Item(154, 17).Value = ""
$ws.Cells.Item(154, 18).Value = ""
$ws.Cells.Item(154, 19).Value = ""
$ws.Cells.Item(154, 20).Value = ""
$ws.Cells.Item(154, 21).Value = ""
$ws.Cells.Item(154, 22).Value = ""
$ws.Cells.Item(154, 23).Value = ""
$ws.Cells.Item(154, 24).Value = ""
$ws.Cells.Item(154, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(154, 26).Value = 10
$ws.Cells.Item(154, 27).Value = $false

# Row 155
$ws.Cells.Item(155, 1).Value = 153
$ws.Cells.Item(155, 2).Value = 52
$ws.Cells.Item(155, 3).Value = "What areas are within 3 minutes of driving time from the nearest fire station (from my current location) in Oleander"
$ws.Cells.Item(155, 4).Value = "Oleander"
$ws.Cells.Item(155, 5).Value = ""
$ws.Cells.Item(155, 6).Value = ""
$ws.Cells.Item(155, 7).Value = ""
$ws.Cells.Item(155, 8).Value = "done"
$ws.Cells.Item(155, 9).Value = "Data queries"
$ws.Cells.Item(155, 10).Value = "Network analysis"
$ws.Cells.Item(155, 11).Value = "Data queries"
$ws.Cells.Item(155, 12).Value = "Network analysis"
$ws.Cells.Item(155, 13).Value = "classification"
$ws.Cells.Item(155, 14).Value = "Data queries"
$ws.Cells.Item(155, 15).Value = "Overlay analysis"
$ws.Cells.Item(155, 16).Value = ""
$ws.Cells.Item(155, 17).Value = ""
$ws.Cells.Item(155, 18).Value = ""
$ws.Cells.Item(155, 19).Value = ""
$ws.Cells.Item(155, 20).Value = ""
$ws.Cells.Item(155, 21).Value = ""
$ws.Cells.Item(155, 22).Value = ""
$ws.Cells.Item(155, 23).Value = "어디로 부터 가장 가까운 소방서"
$ws.Cells.Item(155, 24).Value = ""
$ws.Cells.Item(155, 25).Value = "data queries,network analysis,data queries,network analysis,classification,data queries,overlay analysis,어디로 부터 가장 가까운 소방서"
$ws.Cells.Item(155, 26).Value = 19
$ws.Cells.Item(155, 27).Value = $false

# Row 156
$ws.Cells.Item(156, 1).Value = 154
$ws.Cells.Item(156, 2).Value = 57
$ws.Cells.Item(156, 3).Value = "What areas are within 60 minutes of airports in Crook, Deschutes, and Jefferson county"
$ws.Cells.Item(156, 4).Value = "Crook, Deschutes, Jefferson county"
$ws.Cells.Item(156, 5).Value = ""
$ws.Cells.Item(156, 6).Value = ""
$ws.Cells.Item(156, 7).Value = " aeroway=*"
$ws.Cells.Item(156, 8).Value = "done"
$ws.Cells.Item(156, 9).Value = "Data queries"
$ws.Cells.Item(156, 10).Value = "Network analysis"
$ws.Cells.Item(156, 11).Value = "classification"
$ws.Cells.Item(156, 12).Value = "Data queries"
$ws.Cells.Item(156, 13).Value = "Overlay analysis"
$ws.Cells.Item(156, 14).Value = ""
$ws.Cells.Item(156, 15).Value = ""
$ws.Cells.Item(156, 16).Value = ""
$ws.Cells.Item(156, 17).Value = ""
$ws.Cells.Item(156, 18).Value = ""
$ws.Cells.Item(156, 19).Value = ""
$ws.Cells.Item(156, 20).Value = ""
$ws.Cells.Item(156, 21).Value = ""
$ws.Cells.Item(156, 22).Value = ""
$ws.Cells.Item(156, 23).Value = ""
$ws.Cells.Item(156, 24).Value = ""
$ws.Cells.Item(156, 25).Value = "data queries,network analysis,classification,data queries,overlay analysis"
$ws.Cells.Item(156, 26).Value = 0
$ws.Cells.Item(156, 27).Value = $false

# Row 157
$ws.Cells.Item(157, 1).Value = 155
$ws.Cells.Item(157, 2).Value = 64
$ws.Cells.Item(157, 3).Value = "What areas are within one mile of main roads in Loudoun County in US"
$ws.Cells.Item(157, 4).Value = "Loudoun County"
$ws.Cells.Item(157, 5).Value = "US"
$ws.Cells.Item(157, 6).Value = ""
$ws.Cells.Item(157, 7).Value = "highway=motorway"
$ws.Cells.Item(157, 8).Value = "done"
$ws.Cells.Item(157, 9).Value = "Data queries"
$ws.Cells.Item(157, 10).Value = "Buffer"
$ws.Cells.Item(157, 11).Value = "Overlay analysis"
$ws.Cells.Item(157, 12).Value = ""
$ws.Cells.Item(157, 13).Value = ""
$ws.Cells.Item(157, 14).Value = ""
$ws.Cells.Item(157, 15).Value = ""
$ws.Cells.Item(157, 16).Value = ""
$ws.Cells.Item(157, 17).Value = ""
$ws.Cells.Item(157, 18).Value = ""
$ws.Cells.Item(157, 19).Value = ""
$ws.Cells.Item(157, 20).Value = ""
$ws.Cells.Item(157, 21).Value = ""
$ws.Cells.Item(157, 22).Value = ""
$ws.Cells.Item(157, 23).Value = ""
$ws.Cells.Item(157, 24).Value = ""
$ws.Cells.Item(157, 25).Value = "data queries,buffer,overlay analysis"
$ws.Cells.Item(157, 26).Value = 10
$ws.Cells.Item(157, 27).Value = $false

# Row 158
$ws.Cells.Item(158, 1).Value = 156
$ws.Cells.Item(158, 2).Value = 75
$ws.Cells.Item(158, 3).Value = "What areas have an annual amount of snowfall more than 1000 millimeters in the Cape Peninsula"
$ws.Cells.Item(158, 4).Value = "the Cape Peninsula"
$ws.Cells.Item(158, 5).Value = ""
$ws.Cells.Item(158, 6).Value = ""
$ws.Cells.Item(158, 7).Value = ""
$ws.Cells.Item(158, 8).Value = "done"
$ws.Cells.Item(158, 9).Value = "Geostatistics  "
$ws.Cells.Item(158, 10).Value = "classification"
$ws.Cells.Item(158, 11).Value = "Data queries"
$ws.Cells.Item(158, 12).Value = "Data model conversion"
$ws.Cells.Item(158, 13).Value = "Overlay analysis"
$ws.Cells.Item(158, 14).Value = ""
$ws.Cells.Item(158, 15).Value = ""
$ws.Cells.Item(158, 16).Value = ""
$ws.Cells.Item(158, 17).Value = ""
$ws.Cells.Item(158, 18).Value = ""
$ws.Cells.Item(158, 19).Value = ""
$ws.Cells.Item(158, 20).Value = ""
$ws.Cells.Item(158, 21).Value = ""
$ws.Cells.Item(158, 22).Value = ""
$ws.Cells.Item(158, 23).Value = ""
$ws.Cells.Item(158, 24).Value = ""
$ws.Cells.Item(158, 25).Value = "geostatistics  ,classification,data queries,data model conversion,overlay analysis"
$ws.Cells.Item(158, 26).Value = 26
$ws.Cells.Item(158, 27).Value = $false

# Row 159
$ws.Cells.Item(159, 1).Value = 157
$ws.Cells.Item(159, 2).Value = 77
$ws.Cells.Item(159, 3).Value = "What areas have an annual rainfall of more than 1000 millimeters in the Cape Peninsula"
$ws.Cells.Item(159, 4).Value = "the Cape Peninsula"
$ws.Cells.Item(159, 5).Value = ""
$ws.Cells.Item(159, 6).Value = ""
$ws.Cells.Item(159, 7).Value = ""
$ws.Cells.Item(159, 8).Value = "done"
$ws.Cells.Item(159, 9).Value = "Geostatistics  "
$ws.Cells.Item(159, 10).Value = "classification"
$ws.Cells.Item(159, 11).Value = "Data queries"
$ws.Cells.Item(159, 12).Value = "Data model conversion"
$ws.Cells.Item(159, 13).Value = "Overlay analysis"
$ws.Cells.Item(159, 14).Value = ""
$ws.Cells.Item(159, 15).Value = ""
$ws.Cells.Item(159, 16).Value = ""
$ws.Cells.Item(159, 17).Value = ""
$ws.Cells.Item(159, 18).Value = ""
$ws.Cells.Item(159, 19).Value = ""
$ws.Cells.Item(159, 20).Value = ""
$ws.Cells.Item(159, 21).Value = ""
$ws.Cells.Item(159, 22).Value = ""
$ws.Cells.Item(159, 23).Value = ""
$ws.Cells.Item(159, 24).Value = ""
$ws.Cells.Item(159, 25).Value = "geostatistics  ,classification,data queries,data model conversion,overlay analysis"
$ws.Cells.Item(159, 26).Value = 26
$ws.Cells.Item(159, 27).Value = $false

# Row 160
$ws.Cells.Item(160, 1).Value = 158
$ws.Cells.Item(160, 2).Value = 78
$ws.Cells.Item(160, 3).Value = "What areas have an aspect larger than 45 degree and smaller than 135 degrees in the Cape Peninsula"
$ws.Cells.Item(160, 4).Value = "the Cape Peninsula"
$ws.Cells.Item(160, 5).Value = ""
$ws.Cells.Item(160, 6).Value = ""
$ws.Cells.Item(160, 7).Value = ""
$ws.Cells.Item(160, 8).Value = "done"
$ws.Cells.Item(160, 9).Value = "Topography"
$ws.Cells.Item(160, 10).Value = "classification"
$ws.Cells.Item(160, 11).Value = "Data queries"
$ws.Cells.Item(160, 12).Value = "Data model conversion"
$ws.Cells.Item(160, 13).Value = "Overlay analysis"
$ws.Cells.Item(160, 14).Value = ""
$ws.Cells.Item(160, 15).Value = ""
$ws.Cells.Item(160, 16).Value = ""
$ws.Cells.Item(160, 17).Value = ""
$ws.Cells.Item(160, 18).Value = ""
$ws.Cells.Item(160, 19).Value = ""
$ws.Cells.Item(160, 20).Value = ""
$ws.Cells.Item(160, 21).Value = ""
$ws.Cells.Item(160, 22).Value = ""
$ws.Cells.Item(160, 23).Value = ""
$ws.Cells.Item(160, 24).Value = ""
$ws.Cells.Item(160, 25).Value = "topography,classification,data queries,data model conversion,overlay analysis"
$ws.Cells.Item(160, 26).Value = 24
$ws.Cells.Item(160, 27).Value = $false

# Row 161
$ws.Cells.Item(161, 1).Value = 159
$ws.Cells.Item(161, 2).Value = 83
$ws.Cells.Item(161, 3).Value = "What houses are larger than 30 square meters and within 1km from the nearest school (from my current location) in Utrecht"
$ws.Cells.Item(161, 4).Value = "Utrecht"
$ws.Cells.Item(161, 5).Value = ""
$ws.Cells.Item(161, 6).Value = ""
$ws.Cells.Item(161, 7).Value = "amenity=school, building=house"
$ws.Cells.Item(161, 8).Value = "done"
$ws.Cells.Item(161, 9).Value = "network analysis"
$ws.Cells.Item(161, 10).Value = "Data queries"
$ws.Cells.Item(161, 11).Value = "buffer"
$ws.Cells.Item(161, 12).Value = "Overlay analysis"
$ws.Cells.Item(161, 13).Value = "Data queries"
$ws.Cells.Item(161, 14).Value = ""
$ws.Cells.Item(161, 15).Value = ""
$ws.Cells.Item(161, 16).Value = ""
$ws.Cells.Item(161, 17).Value = ""
$ws.Cells.Item(161, 18).Value = ""
$ws.Cells.Item(161, 19).Value = ""
$ws.Cells.Item(161, 20).Value = ""
$ws.Cells.Item(161, 21).Value = ""
$ws.Cells.Item(161, 22).Value = ""
$ws.Cells.Item(161, 23).Value = ""
$ws.Cells.Item(161, 24).Value = ""
$ws.Cells.Item(161, 25).Value = "network analysis,data queries,buffer,overlay analysis,data queries"
$ws.Cells.Item(161, 26).Value = 29
$ws.Cells.Item(161, 27).Value = $false

# Row 162
$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(162, 2).Value = 86
$ws.Cells.Item(162, 3).Value = "What houses are for sale and within 1km from the nearest school (from my current location) in Utrecht"
$ws.Cells.Item(162, 4).Value = "Utrecht"
$ws.Cells.Item(162, 5).Value = ""
$ws.Cells.Item(162, 6).Value = ""
$ws.Cells.Item(162, 7).Value = "amenity=school, building=house"
$ws.Cells.Item(162, 8).Value = "done"
$ws.Cells.Item(162, 9).Value = "network analysis"
$ws.Cells.Item(162, 10).Value = "Data queries"
$ws.Cells.Item(162, 11).Value = "buffer"
$ws.Cells.Item(162, 12).Value = "Overlay analysis"
$ws.Cells.Item(162, 13).Value = "Data queries"
$ws.Cells.Item(162, 14).Value = ""
$ws.Cells.Item(162, 15).Value = ""
$ws.Cells.Item(162, 16).Value = ""
$ws.Cells.Item(162, 17).Value = ""
$ws.Cells.Item(162, 18).Value = ""
$ws.Cells.Item(162, 19).Value = ""
$ws.Cells.Item(162, 20).Value = ""
$ws.Cells.Item(162, 21).Value = ""
$ws.Cells.Item(162, 22).Value = ""
$ws.Cells.Item(162, 23).Value = ""
$ws.Cells.Item(162, 24).Value = ""
$ws.Cells.Item(162, 25).Value = "network analysis,data queries,buffer,overlay analysis,data queries"
$ws.Cells.Item(162, 26).Value = 29
$ws.Cells.Item(162, 27).Value = $false

# Row 163
$ws.Cells.Item(163, 1).Value = 161
$ws.Cells.Item(163, 2).Value = 90
$ws.Cells.Item(163, 3).Value = "What houses are for sale in flood zone in Utrecht"
$ws.Cells.Item(163, 4).Value = "Utrecht"
$ws.Cells.Item(163, 5).Value = ""
$ws.Cells.Item(163, 6).Value = ""
$ws.Cells.Item(163, 7).Value = "boundary=urban, building=house"
$ws.Cells.Item(163, 8).Value = "done"
$ws.Cells.Item(163, 9).Value = "Data queries"
$ws.Cells.Item(163, 10).Value = "Overlay analysis"
$ws.Cells.Item(163, 11).Value = "Data queries"
$ws.Cells.Item(163, 12).Value = ""
$ws.Cells.Item(163, 13).Value = ""
$ws.Cells.Item(163, 14).Value = ""
$ws.Cells.Item(163, 15).Value = ""
$ws.Cells.Item(163, 16).Value = ""
$ws.Cells.Item(163, 17).Value = ""
$ws.Cells.Item(163, 18).Value = ""
$ws.Cells.Item(163, 19).Value = ""
$ws.Cells.Item(163, 20).Value = ""
$ws.Cells.Item(163, 21).Value = ""
$ws.Cells.Item(163, 22).Value = ""
$ws.Cells.Item(163, 23).Value = "osm urban이 있긴 하지만 거의 안쓴다"
$ws.Cells.Item(163, 24).Value = ""
$ws.Cells.Item(163, 25).Value = "data queries,overlay analysis,data queries,osm urban이 있긴 하지만 거의 안쓴다"
$ws.Cells.Item(163, 26).Value = 30
$ws.Cells.Item(163, 27).Value = $false

# Row 164
$ws.Cells.Item(164, 1).Value = 162
$ws.Cells.Item(164, 2).Value = 93
$ws.Cells.Item(164, 3).Value = "What houses are greather than 30 square meters in urban areas in Utrecht"
$ws.Cells.Item(164, 4).Value = "Utrecht"
$ws.Cells.Item(164, 5).Value = ""
$ws.Cells.Item(164, 6).Value = ""
$ws.Cells.Item(164, 7).Value = "boundary=urban, building=house"
$ws.Cells.Item(164, 8).Value = "done"
$ws.Cells.Item(164, 9).Value = "Data queries"
$ws.Cells.Item(164, 10).Value = "Overlay analysis"
$ws.Cells.Item(164, 11).Value = "Data queries"
$ws.Cells.Item(164, 12).Value = ""
$ws.Cells.Item(164, 13).Value = ""
$ws.Cells.Item(164, 14).Value = ""
$ws.Cells.Item(164, 15).Value = ""
$ws.Cells.Item(164, 16).Value = ""
$ws.Cells.Item(164, 17).Value = ""
$ws.Cells.Item(164, 18).Value = ""
$ws.Cells.Item(164, 19).Value = ""
$ws.Cells.Item(164, 20).Value = ""
$ws.Cells.Item(164, 21).Value = ""
$ws.Cells.Item(164, 22).Value = ""
$ws.Cells.Item(164, 23).Value = "osm urban이 있긴 하지만 거의 안쓴다"
$ws.Cells.Item(164, 24).Value = ""
$ws.Cells.Item(164, 25).Value = "data queries,overlay analysis,data queries,osm urban이 있긴 하지만 거의 안쓴다"
$ws.Cells.Item(164, 26).Value = 30
$ws.Cells.Item(164, 27).Value = $false

# Row 165
$ws.Cells.Item(165, 1).Value = 163
$ws.Cells.Item(165, 2).Value = 105
$ws.Cells.Item(165, 3).Value = "What is the central feature of bank branches in Oleander"
$ws.Cells.Item(165, 4).Value = "Oleander"
$ws.Cells.Item(165, 5).Value = ""
$ws.Cells.Item(165, 6).Value = ""
$ws.Cells.Item(165, 7).Value = "amenity=bank"
$ws.Cells.Item(165, 8).Value = "done"
$ws.Cells.Item(165, 9).Value = "Data queries"
$ws.Cells.Item(165, 10).Value = "Generalization"
$ws.Cells.Item(165, 11).Value = "Geostatistics  "
$ws.Cells.Item(165, 12).Value = ""
$ws.Cells.Item(165, 13).Value = ""
$ws.Cells.Item(165, 14).Value = ""
$ws.Cells.Item(165, 15).Value = ""
$ws.Cells.Item(165, 16).Value = ""
$ws.Cells.Item(165, 17).Value = ""
$ws.Cells.Item(165, 18).Value = ""
$ws.Cells.Item(165, 19).Value = ""
$ws.Cells.Item(165, 20).Value = ""
$ws.Cells.Item(165, 21).Value = ""
$ws.Cells.Item(165, 22).Value = "https://pro.arcgis.com/en/pro-app/latest/tool-reference/spatial-statistics/h-how-central-feature-spatial-statistics-works.htm"
$ws.Cells.Item(165, 23).Value = ""
$ws.Cells.Item(165, 24).Value = ""
$ws.Cells.Item(165, 25).Value = "data queries,generalization,geostatistics  ,https://pro.arcgis.com/en/pro-app/latest/tool-reference/spatial-statistics/h-how-central-feature-spatial-statistics-works.htm"
$ws.Cells.Item(165, 26).Value = 36
$ws.Cells.Item(165, 27).Value = $false

# Row 166
$ws.Cells.Item(166, 1).Value = 164
$ws.Cells.Item(166, 2).Value = 112
$ws.Cells.Item(166, 3).Value = "What is the density surface of temperature measurements in Oleander city"
$ws.Cells.Item(166, 4).Value = "Oleander city"
$ws.Cells.Item(166, 5).Value = ""
$ws.Cells.Item(166, 6).Value = ""
$ws.Cells.Item(166, 7).Value = " man_made=monitoring_station"
$ws.Cells.Item(166, 8).Value = "done"
$ws.Cells.Item(166, 9).Value = "Data queries"
$ws.Cells.Item(166, 10).Value = "Geostatistics  "
$ws.Cells.Item(166, 11).Value = ""
$ws.Cells.Item(166, 12).Value = ""
$ws.Cells.Item(166, 13).Value = ""
$ws.Cells.Item(166, 14).Value = ""
$ws.Cells.Item(166, 15).Value = ""
$ws.Cells.Item(166, 16).Value = ""
$ws.Cells.Item(166, 17).Value = ""
$ws.Cells.Item(166, 18).Value = ""
$ws.Cells.Item(166, 19).Value = ""
$ws.Cells.Item(166, 20).Value = ""
$ws.Cells.Item(166, 21).Value = ""
$ws.Cells.Item(166, 22).Value = ""
$ws.Cells.Item(166, 23).Value = ""
$ws.Cells.Item(166, 24).Value = ""
$ws.Cells.Item(166, 25).Value = "data queries,geostatistics  "
$ws.Cells.Item(166, 26).Value = 41
$ws.Cells.Item(166, 27).Value = $false

# Row 167
$ws.Cells.Item(167, 1).Value = 165
$ws.Cells.Item(167, 2).Value = 115
$ws.Cells.Item(167, 3).Value = "What is the Euclidean distance to green areas in Amsterdam"
$ws.Cells.Item(167, 4).Value = "Amsterdam"
$ws.Cells.Item(167, 5).Value = ""
$ws.Cells.Item(167, 6).Value = ""
$ws.Cells.Item(167, 7).Value = ""
$ws.Cells.Item(167, 8).Value = "done"
$ws.Cells.Item(167, 9).Value = "Data queries"
$ws.Cells.Item(167, 10).Value = "Geostatistics  "
$ws.Cells.Item(167, 11).Value = ""
$ws.Cells.Item(167, 12).Value = ""
$ws.Cells.Item(167, 13).Value = ""
$ws.Cells.Item(167, 14).Value = ""
$ws.Cells.Item(167, 15).Value = ""
$ws.Cells.Item(167, 16).Value = ""
$ws.Cells.Item(167, 17).Value = ""
$ws.Cells.Item(167, 18).Value = ""
$ws.Cells.Item(167, 19).Value = ""
$ws.Cells.Item(167, 20).Value = ""
$ws.Cells.Item(167, 21).Value = ""
$ws.Cells.Item(167, 22).Value = ""
$ws.Cells.Item(167, 23).Value = ""
$ws.Cells.Item(167, 24).Value = ""
$ws.Cells.Item(167, 25).Value = "data queries,geostatistics  "
$ws.Cells.Item(167, 26).Value = 41
$ws.Cells.Item(167, 27).Value = $false

# Row 168
$ws.Cells.Item(168, 1).Value = 166
$ws.Cells.Item(168, 2).Value = 120
$ws.Cells.Item(168, 3).Value = "What is the Euclidean distance to the rivers in Crook, Deschutes, and Jefferson county"
$ws.Cells.Item(168, 4).Value = "Crook, Deschutes, Jefferson county"
$ws.Cells.Item(168, 5).Value = ""
$ws.Cells.Item(168, 6).Value = ""
$ws.Cells.Item(168, 7).Value = ""
$ws.Cells.Item(168, 8).Value = "done"
$ws.Cells.Item(168, 9).Value = "Data queries"
$ws.Cells.Item(168, 10).Value = "Geostatistics  "
$ws.Cells.Item(168, 11).Value = ""
$ws.Cells.Item(168, 12).Value = ""
$ws.Cells.Item(168, 13).Value = ""
$ws.Cells.Item(168, 14).Value = ""
$ws.Cells.Item(168, 15).Value = ""
$ws.Cells.Item(168, 16).Value = ""
$ws.Cells.Item(168, 17).Value = ""
$ws.Cells.Item(168, 18).Value = ""
$ws.Cells.Item(168, 19).Value = ""
$ws.Cells.Item(168, 20).Value = ""
$ws.Cells.Item(168, 21).Value = ""
$ws.Cells.Item(168, 22).Value = ""
$ws.Cells.Item(168, 23).Value = ""
$ws.Cells.Item(168, 24).Value = ""
$ws.Cells.Item(168, 25).Value = "data queries,geostatistics  "
$ws.Cells.Item(168, 26).Value = 41
$ws.Cells.Item(168, 27).Value = $false

# Row 169
$ws.Cells.Item(169, 1).Value = 167
$ws.Cells.Item(169, 2).Value = 127
$ws.Cells.Item(169, 3).Value = "What is the lung cancer mortality rate of white males for each city in the Western USA from 1970 to 1994"
$ws.Cells.Item(169, 4).Value = "the Western USA"
$ws.Cells.Item(169, 5).Value = ""
$ws.Cells.Item(169, 6).Value = " from 1970 to 1994"
$ws.Cells.Item(169, 7).Value = ""
$ws.Cells.Item(169, 8).Value = "done"
$ws.Cells.Item(169, 9).Value = "data editing"
$ws.Cells.Item(169, 10).Value = "Data queries"
$ws.Cells.Item(169, 11).Value = ""
$ws.Cells.Item(169, 12).Value = ""
$ws.Cells.Item(169, 13).Value = ""
$ws.Cells.Item(169, 14).Value = ""
$ws.Cells.Item(169, 15).Value = ""
$ws.Cells.Item(169, 16).Value = ""
$ws.Cells.Item(169, 17).Value = ""
$ws.Cells.Item(169, 18).Value = ""
$ws.Cells.Item(169, 19).Value = ""
$ws.Cells.Item(169, 20).Value = ""
$ws.Cells.Item(169, 21).Value = ""
$ws.Cells.Item(169, 22).Value = ""
$ws.Cells.Item(169, 23).Value = ""
$ws.Cells.Item(169, 24).Value = ""
$ws.Cells.Item(169, 25).Value = "data editing,data queries"
$ws.Cells.Item(169, 26).Value = 9
$ws.Cells.Item(169, 27).Value = $false

# Row 170
$ws.Cells.Item(170, 1).Value = 168
$ws.Cells.Item(170, 2).Value = 129
$ws.Cells.Item(170, 3).Value = "What is the mean center of customers weighted by the number of transactions in Oleander city"
$ws.Cells.Item(170, 4).Value = "Oleander city"
$ws.Cells.Item(170, 5).Value = ""
$ws.Cells.Item(170, 6).Value = ""
$ws.Cells.Item(170, 7).Value = ""
$ws.Cells.Item(170, 8).Value = "done"
$ws.Cells.Item(170, 9).Value = "Data queries"
$ws.Cells.Item(170, 10).Value = "Geostatistics  "
$ws.Cells.Item(170, 11).Value = ""
$ws.Cells.Item(170, 12).Value = ""
$ws.Cells.Item(170, 13).Value = ""
$ws.Cells.Item(170, 14).Value = ""
$ws.Cells.Item(170, 15).Value = ""
$ws.Cells.Item(170, 16).Value = ""
$ws.Cells.Item(170, 17).Value = ""
$ws.Cells.Item(170, 18).Value = ""
$ws.Cells.Item(170, 19).Value = ""
$ws.Cells.Item(170, 20).Value = ""
$ws.Cells.Item(170, 21).Value = ""
$ws.Cells.Item(170, 22).Value = ""
$ws.Cells.Item(170, 23).Value = ""
$ws.Cells.Item(170, 24).Value = ""
$ws.Cells.Item(170, 25).Value = "data queries,geostatistics  "
$ws.Cells.Item(170, 26).Value = 41
$ws.Cells.Item(170, 27).Value = $false

# Row 171
$ws.Cells.Item(171, 1).Value = 169
$ws.Cells.Item(171, 2).Value = 134
$ws.Cells.Item(171, 3).Value = "What is the mean center of the fire calls weighted by the priority in Fort Worth"
$ws.Cells.Item(171, 4).Value = "Fort Worth"
$ws.Cells.Item(171, 5).Value = ""
$ws.Cells.Item(171, 6).Value = ""
$ws.Cells.Item(171, 7).Value = ""
$ws.Cells.Item(171, 8).Value = "done"
$ws.Cells.Item(171, 9).Value = "Data queries"
$ws.Cells.Item(171, 10).Value = "Overlay analysis"
$ws.Cells.Item(171, 11).Value = "Geostatistics  "
$ws.Cells.Item(171, 12).Value = ""
$ws.Cells.Item(171, 13).Value = ""
$ws.Cells.Item(171, 14).Value = ""
$ws.Cells.Item(171, 15).Value = ""
$ws.Cells.Item(171, 16).Value = ""
$ws.Cells.Item(171, 17).Value = ""
$ws.Cells.Item(171, 18).Value = ""
$ws.Cells.Item(171, 19).Value = ""
$ws.Cells.Item(171, 20).Value = ""
$ws.Cells.Item(171, 21).Value = ""
$ws.Cells.Item(171, 22).Value = ""
$ws.Cells.Item(171, 23).Value = ""
$ws.Cells.Item(171, 24).Value = ""
$ws.Cells.Item(171, 25).Value = "data queries,overlay analysis,geostatistics  "
$ws.Cells.Item(171, 26).Value = 44
$ws.Cells.Item(171, 27).Value = $false

# Row 172
$ws.Cells.Item(172, 1).Value = 170
$ws.Cells.Item(172, 2).Value = 139
$ws.Cells.Item(172, 3).Value = "What is the median household income for each census block in Tarrant County in Texas"
$ws.Cells.Item(172, 4).Value = "Tarrant County"
$ws.Cells.Item(172, 5).Value = " Texas"
$ws.Cells.Item(172, 6).Value = ""
$ws.Cells.Item(172, 7).Value = ""
$ws.Cells.Item(172, 8).Value = "done"
$ws.Cells.Item(172, 9).Value = "data editing"
$ws.Cells.Item(172, 10).Value = "Overlay analysis"
$ws.Cells.Item(172, 11).Value = "data editing"
$ws.Cells.Item(172, 12).Value = "Data queries"
$ws.Cells.Item(172, 13).Value = ""
$ws.Cells.Item(172, 14).Value = ""
$ws.Cells.Item(172, 15).Value = ""
$ws.Cells.Item(172, 16).Value = ""
$ws.Cells.Item(172, 17).Value = ""
$ws.Cells.Item(172, 18).Value = ""
$ws.Cells.Item(172, 19).Value = ""
$ws.Cells.Item(172, 20).Value = ""
$ws.Cells.Item(172, 21).Value = ""
$ws.Cells.Item(172, 22).Value = ""
$ws.Cells.Item(172, 23).Value = ""
$ws.Cells.Item(172, 24).Value = ""
$ws.Cells.Item(172, 25).Value = "data editing,overlay analysis,data editing,data queries"
$ws.Cells.Item(172, 26).Value = 34
$ws.Cells.Item(172, 27).Value = $false

# Row 173
$ws.Cells.Item(173, 1).Value = 171
$ws.Cells.Item(173, 2).Value = 141
$ws.Cells.Item(173, 3).Value = "What is the median people age for each census tract in Tarrant County in Texas"
$ws.Cells.Item(173, 4).Value = "Tarrant County"
$ws.Cells.Item(173, 5).Value = " Texas"
$ws.Cells.Item(173, 6).Value = ""
$ws.Cells.Item(173, 7).Value = ""
$ws.Cells.Item(173, 8).Value = "done"
$ws.Cells.Item(173, 9).Value = "data editing"
$ws.Cells.Item(173, 10).Value = "Overlay analysis"
$ws.Cells.Item(173, 11).Value = "data editing"
$ws.Cells.Item(173, 12).Value = "Data queries"
$ws.Cells.Item(173, 13).Value = ""
$ws.Cells.Item(173, 14).Value = ""
$ws.Cells.Item(173, 15).Value = ""
$ws.Cells.Item(173, 16).Value = ""
$ws.Cells.Item(173, 17).Value = ""
$ws.Cells.Item(173, 18).Value = ""
$ws.Cells.Item(173, 19).Value = ""
$ws.Cells.Item(173, 20).Value = ""
$ws.Cells.Item(173, 21).Value = ""
$ws.Cells.Item(173, 22).Value = ""
$ws.Cells.Item(173, 23).Value = ""
$ws.Cells.Item(173, 24).Value = ""
$ws.Cells.Item(173, 25).Value = "data editing,overlay analysis,data editing,data queries"
$ws.Cells.Item(173, 26).Value = 34
$ws.Cells.Item(173, 27).Value = $false

# Row 174
$ws.Cells.Item(174, 1).Value = 172
$ws.Cells.Item(174, 2).Value = 149
$ws.Cells.Item(174, 3).Value = "What liquor stores are within 1000 foot of schools in El Cajon"
$ws.Cells.Item(174, 4).Value = "El Cajon"
$ws.Cells.Item(174, 5).Value = ""
$ws.Cells.Item(174, 6).Value = ""
$ws.Cells.Item(174, 7).Value = "shop=alcohol, amenity=school"
$ws.Cells.Item(174, 8).Value = "done"
$ws.Cells.Item(174, 9).Value = "Data queries"
$ws.Cells.Item(174, 10).Value = "buffer"
$ws.Cells.Item(174, 11).Value = "Overlay analysis"
$ws.Cells.Item(174, 12).Value = "Data queries"
$ws.Cells.Item(174, 13).Value = ""
$ws.Cells.Item(174, 14).Value = ""
$ws.Cells.Item(174, 15).Value = ""
$ws.Cells.Item(174, 16).Value = ""
$ws.Cells.Item(174, 17).Value = ""
$ws.Cells.Item(174, 18).Value = ""
$ws.Cells.Item(174, 19).Value = ""
$ws.Cells.Item(174, 20).Value = ""
$ws.Cells.Item(174, 21).Value = ""
$ws.Cells.Item(174, 22).Value = ""
$ws.Cells.Item(174, 23).Value = ""
$ws.Cells.Item(174, 24).Value = ""
$ws.Cells.Item(174, 25).Value = "data queries,buffer,overlay analysis,data queries"
$ws.Cells.Item(174, 26).Value = 28
$ws.Cells.Item(174, 27).Value = $false

# Row 175
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = 153
$ws.Cells.Item(175, 3).Value = "Where are not protected region in Assam in India"
$ws.Cells.Item(175, 4).Value = "Assam"
$ws.Cells.Item(175, 5).Value = "India"
$ws.Cells.Item(175, 6).Value = ""
$ws.Cells.Item(175, 7).Value = "landuse=conservation"
$ws.Cells.Item(175, 8).Value = "done"
$ws.Cells.Item(175, 9).Value = "Data queries"
$ws.Cells.Item(175, 10).Value = "Geometry measurement"
$ws.Cells.Item(175, 11).Value = "Data queries"
$ws.Cells.Item(175, 12).Value = ""
$ws.Cells.Item(175, 13).Value = ""
$ws.Cells.Item(175, 14).Value = ""
$ws.Cells.Item(175, 15).Value = ""
$ws.Cells.Item(175, 16).Value = ""
$ws.Cells.Item(175, 17).Value = ""
$ws.Cells.Item(175, 18).Value = ""
$ws.Cells.Item(175, 19).Value = ""
$ws.Cells.Item(175, 20).Value = ""
$ws.Cells.Item(175, 21).Value = ""
$ws.Cells.Item(175, 22).Value = ""
$ws.Cells.Item(175, 23).Value = ""
$ws.Cells.Item(175, 24).Value = ""
$ws.Cells.Item(175, 25).Value = "data queries,geometry measurement,data queries"
$ws.Cells.Item(175, 26).Value = 47
$ws.Cells.Item(175, 27).Value = $false

# Row 176
$ws.Cells.Item(176, 1).Value = 174
$ws.Cells.Item(176, 2).Value = 156
$ws.Cells.Item(176, 3).Value = "Where are the commercial areas in Amsterdam"
$ws.Cells.Item(176, 4).Value = "Amsterdam"
$ws.Cells.Item(176, 5).Value = ""
$ws.Cells.Item(176, 6).Value = ""
$ws.Cells.Item(176, 7).Value = "landuse=commercial"
$ws.Cells.Item(176, 8).Value = "done"
$ws.Cells.Item(176, 9).Value = "Data queries"
$ws.Cells.Item(176, 10).Value = "Geometry measurement"
$ws.Cells.Item(176, 11).Value = "Data queries"
$ws.Cells.Item(176, 12).Value = ""
$ws.Cells.Item(176, 13).Value = ""
$ws.Cells.Item(176, 14).Value = ""
$ws.Cells.Item(176, 15).Value = ""
$ws.Cells.Item(176, 16).Value = ""
$ws.Cells.Item(176, 17).Value = ""
$ws.Cells.Item(176, 18).Value = ""
$ws.Cells.Item(176, 19).Value = ""
$ws.Cells.Item(176, 20).Value = ""
$ws.Cells.Item(176, 21).Value = ""
$ws.Cells.Item(176, 22).Value = ""
$ws.Cells.Item(176, 23).Value = ""
$ws.Cells.Item(176, 24).Value = ""
$ws.Cells.Item(176, 25).Value = "data queries,geometry measurement,data queries"
$ws.Cells.Item(176, 26).Value = 47
$ws.Cells.Item(176, 27).Value = $false

# Row 177
$ws.Cells.Item(177, 1).Value = 175
$ws.Cells.Item(177, 2).Value = 164
$ws.Cells.Item(177, 3).Value = "Where are the rocky areas in Spain"
$ws.Cells.Item(177, 4).Value = "Spain"
$ws.Cells.Item(177, 5).Value = ""
$ws.Cells.Item(177, 6).Value = ""
$ws.Cells.Item(177, 7).Value = "nature : bare_rock"
$ws.Cells.Item(177, 8).Value = "done"
$ws.Cells.Item(177, 9).Value = "Data queries"
$ws.Cells.Item(177, 10).Value = "Geometry measurement"
$ws.Cells.Item(177, 11).Value = "Data queries"
$ws.Cells.Item(177, 12).Value = ""
$ws.Cells.Item(177, 13).Value = ""
$ws.Cells.Item(177, 14).Value = ""
$ws.Cells.Item(177, 15).Value = ""
$ws.Cells.Item(177, 16).Value = ""
$ws.Cells.Item(177, 17).Value = ""
$ws.Cells.Item(177, 18).Value = ""
$ws.Cells.Item(177, 19).Value = ""
$ws.Cells.Item(177, 20).Value = ""
$ws.Cells.Item(177, 21).Value = ""
$ws.Cells.Item(177, 22).Value = ""
$ws.Cells.Item(177, 23).Value = ""
$ws.Cells.Item(177, 24).Value = ""
$ws.Cells.Item(177, 25).Value = "data queries,geometry measurement,data queries"
$ws.Cells.Item(177, 26).Value = 47
$ws.Cells.Item(177, 27).Value = $false

# Row 178
$ws.Cells.Item(178, 1).Value = 176
$ws.Cells.Item(178, 2).Value = 167
$ws.Cells.Item(178, 3).Value = "Which houses are for sale in Utrecht"
$ws.Cells.Item(178, 4).Value = "Utrecht"
$ws.Cells.Item(178, 5).Value = ""
$ws.Cells.Item(178, 6).Value = ""
$ws.Cells.Item(178, 7).Value = "building=house"
$ws.Cells.Item(178, 8).Value = "done"
$ws.Cells.Item(178, 9).Value = "data editing"
$ws.Cells.Item(178, 10).Value = "Data queries"
$ws.Cells.Item(178, 11).Value = ""
$ws.Cells.Item(178, 12).Value = ""
$ws.Cells.Item(178, 13).Value = ""
$ws.Cells.Item(178, 14).Value = ""
$ws.Cells.Item(178, 15).Value = ""
$ws.Cells.Item(178, 16).Value = ""
$ws.Cells.Item(178, 17).Value = ""
$ws.Cells.Item(178, 18).Value = ""
$ws.Cells.Item(178, 19).Value = ""
$ws.Cells.Item(178, 20).Value = ""
$ws.Cells.Item(178, 21).Value = ""
$ws.Cells.Item(178, 22).Value = ""
$ws.Cells.Item(178, 23).Value = ""
$ws.Cells.Item(178, 24).Value = ""
$ws.Cells.Item(178, 25).Value = "data editing,data queries"
$ws.Cells.Item(178, 26).Value = 9
$ws.Cells.Item(178, 27).Value = $false

# Row 179
$ws.Cells.Item(179, 1).Value = 177
$ws.Cells.Item(179, 2).Value = 170
$ws.Cells.Item(179, 3).Value = "Which houses have construction year between 1990 and 2000 in Utrecht"
$ws.Cells.Item(179, 4).Value = "Utrecht"
$ws.Cells.Item(179, 5).Value = ""
$ws.Cells.Item(179, 6).Value = ""
$ws.Cells.Item(179, 7).Value = "year_of_construction=*"
$ws.Cells.Item(179, 8).Value = "done"
$ws.Cells.Item(179, 9).Value = "Data queries"
$ws.Cells.Item(179, 10).Value = ""
$ws.Cells.Item(179, 11).Value = ""
$ws.Cells.Item(179, 12).Value = ""
$ws.Cells.Item(179, 13).Value = ""
$ws.Cells.Item(179, 14).Value = ""
$ws.Cells.Item(179, 15).Value = ""
$ws.Cells.Item(179, 16).Value = ""
$ws.Cells.Item(179, 17).Value = ""
$ws.Cells.Item(179, 18).Value = ""
$ws.Cells.Item(179, 19).Value = ""
$ws.Cells.Item(179, 20).Value = ""
$ws.Cells.Item(179, 21).Value = ""
$ws.Cells.Item(179, 22).Value = ""
$ws.Cells.Item(179, 23).Value = ""
$ws.Cells.Item(179, 24).Value = ""
$ws.Cells.Item(179, 25).Value = "data queries"
$ws.Cells.Item(179, 26).Value = 1
$ws.Cells.Item(179, 27).Value = $false

# Row 180
$ws.Cells.Item(180, 1).Value = 178
$ws.Cells.Item(180, 2).Value = 177
$ws.Cells.Item(180, 3).Value = "Which schools are not within 3 minutes of driving time from a fire station in Fort Worth"
$ws.Cells.Item(180, 4).Value = "Fort Worth"
$ws.Cells.Item(180, 5).Value = ""
$ws.Cells.Item(180, 6).Value = ""
$ws.Cells.Item(180, 7).Value = "amenity=fire_station, amenity=school"
$ws.Cells.Item(180, 8).Value = "done"
$ws.Cells.Item(180, 9).Value = "Data queries"
$ws.Cells.Item(180, 10).Value = "Network analysis"
$ws.Cells.Item(180, 11).Value = "classification"
$ws.Cells.Item(180, 12).Value = "Data queries"
$ws.Cells.Item(180, 13).Value = "Overlay analysis"
$ws.Cells.Item(180, 14).Value = "data queries"
$ws.Cells.Item(180, 15).Value = ""
$ws.Cells.Item(180, 16).Value = ""
$ws.Cells.Item(180, 17).Value = ""
$ws.Cells.Item(180, 18).Value = ""
$ws.Cells.Item(180, 19).Value = ""
$ws.Cells.Item(180, 20).Value = ""
$ws.Cells.Item(180, 21).Value = ""
$ws.Cells.Item(180, 22).Value = ""
$ws.Cells.Item(180, 23).Value = ""
$ws.Cells.Item(180, 24).Value = ""
$ws.Cells.Item(180, 25).Value = "data queries,network analysis,classification,data queries,overlay analysis,data queries"
$ws.Cells.Item(180, 26).Value = 5
$ws.Cells.Item(180, 27).Value = $false

# Row 181
$ws.Cells.Item(181, 1).Value = 179
$ws.Cells.Item(181, 2).Value = 179
$ws.Cells.Item(181, 3).Value = "Which vacant lots are within 1 mile of a freeway in Hillsboro"
$ws.Cells.Item(181, 4).Value = "Hillsboro"
$ws.Cells.Item(181, 5).Value = ""
$ws.Cells.Item(181, 6).Value = ""
$ws.Cells.Item(181, 7).Value = "abandoned:*=*, highway=motorway"
$ws.Cells.Item(181, 8).Value = "done"
$ws.Cells.Item(181, 9).Value = "Data queries"
$ws.Cells.Item(181, 10).Value = "buffer"
$ws.Cells.Item(181, 11).Value = "Overlay analysis"
$ws.Cells.Item(181, 12).Value = "Data queries"
$ws.Cells.Item(181, 13).Value = ""
$ws.Cells.Item(181, 14).Value = ""
$ws.Cells.Item(181, 15).Value = ""
$ws.Cells.Item(181, 16).Value = ""
$ws.Cells.Item(181, 17).Value = ""
$ws.Cells.Item(181, 18).Value = ""
$ws.Cells.Item(181, 19).Value = ""
$ws.Cells.Item(181, 20).Value = ""
$ws.Cells.Item(181, 21).Value = ""
$ws.Cells.Item(181, 22).Value = ""
$ws.Cells.Item(181, 23).Value = ""
$ws.Cells.Item(181, 24).Value = ""
$ws.Cells.Item(181, 25).Value = "data queries,buffer,overlay analysis,data queries"
$ws.Cells.Item(181, 26).Value = 28
$ws.Cells.Item(181, 27).Value = $false

# Row 182
$ws.Cells.Item(182, 1).Value = 180
$ws.Cells.Item(182, 2).Value = 183
$ws.Cells.Item(182, 3).Value = "Which wind farm proposals are nearest to the high way in Scotland"
$ws.Cells.Item(182, 4).Value = "Scotland"
$ws.Cells.Item(182, 5).Value = ""
$ws.Cells.Item(182, 6).Value = ""
$ws.Cells.Item(182, 7).Value = "highway=*"
$ws.Cells.Item(182, 8).Value = "done"
$ws.Cells.Item(182, 9).Value = "Data queries"
$ws.Cells.Item(182, 10).Value = "network analysis"
$ws.Cells.Item(182, 11).Value = "Data queries"
$ws.Cells.Item(182, 12).Value = ""
$ws.Cells.Item(182, 13).Value = ""
$ws.Cells.Item(182, 14).Value = ""
$ws.Cells.Item(182, 15).Value = ""
$ws.Cells.Item(182, 16).Value = ""
$ws.Cells.Item(182, 17).Value = ""
$ws.Cells.Item(182, 18).Value = ""
$ws.Cells.Item(182, 19).Value = ""
$ws.Cells.Item(182, 20).Value = ""
$ws.Cells.Item(182, 21).Value = ""
$ws.Cells.Item(182, 22).Value = ""
$ws.Cells.Item(182, 23).Value = ""
$ws.Cells.Item(182, 24).Value = ""
$ws.Cells.Item(182, 25).Value = "data queries,network analysis,data queries"
$ws.Cells.Item(182, 26).Value = 4
$ws.Cells.Item(182, 27).Value = $false

# Row 183
$ws.Cells.Item(183, 1).Value = 181
$ws.Cells.Item(183, 2).Value = 185
$ws.Cells.Item(183, 3).Value = "Which wind farm proposals are nearest to the roads in Scotland"
$ws.Cells.Item(183, 4).Value = "Scotland"
$ws.Cells.Item(183, 5).Value = ""
$ws.Cells.Item(183, 6).Value = ""
$ws.Cells.Item(183, 7).Value = "highway=*"
$ws.Cells.Item(183, 8).Value = "done"
$ws.Cells.Item(183, 9).Value = "Data queries"
$ws.Cells.Item(183, 10).Value = "network analysis"
$ws.Cells.Item(183, 11).Value = "Data queries"
$ws.Cells.Item(183, 12).Value = ""
$ws.Cells.Item(183, 13).Value = ""
$ws.Cells.Item(183, 14).Value = ""
$ws.Cells.Item(183, 15).Value = ""
$ws.Cells.Item(183, 16).Value = ""
$ws.Cells.Item(183, 17).Value = ""
$ws.Cells.Item(183, 18).Value = ""
$ws.Cells.Item(183, 19).Value = ""
$ws.Cells.Item(183, 20).Value = ""
$ws.Cells.Item(183, 21).Value = ""
$ws.Cells.Item(183, 22).Value = ""
$ws.Cells.Item(183, 23).Value = ""
$ws.Cells.Item(183, 24).Value = ""
$ws.Cells.Item(183, 25).Value = "data queries,network analysis,data queries"
$ws.Cells.Item(183, 26).Value = 4
$ws.Cells.Item(183, 27).Value = $false
